$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.3715

$ws.Range("C4").Value = -13.9227
$ws.Range("D4").Value = -8.454800000000004
$ws.Range("E4").Value = 11.9415

$ws.Range("D5").Value = -8.303400000000003

$ws.Range("C6").Value = -11.3982
$ws.Range("D6").Value = -8.879599999999989

$ws.Range("C7").Value = -11.5096

$ws.Range("C8").Value = -12.4563
$ws.Range("D8").Value = -8.479699999999999

$ws.Range("E9").Value = 13.37520000000001

$ws.Range("E11").Value = 13.3117

$ws.Range("E14").Value = 13.20270000000001

$ws.Range("C16").Value = -11.8376
$ws.Range("D16").Value = -8.736000000000004

$ws.Range("E18").Value = 12.4804

$ws.Range("C20").Value = -14.50039999999999

$ws.Range("C21").Value = -12.83030000000001

$ws.Range("D22").Value = -7.848799999999998

$ws.Range("E25").Value = 13.07719999999999
